$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H42").Value = 109.875
$ws.Range("I42").Value = 109.875
$ws.Range("K42").Value = 329.625
$ws.Range("M42").Value = -99.625
$ws.Range("H51").Value = 2221
$ws.Range("J51").Value = 2221
$ws.Range("L51").Value = 2221
$ws.Range("N51").Value = -3189
$ws.Range("H115").Value = 682.2222
$ws.Range("I115").Value = 682.2222
$ws.Range("K115").Value = 2046.6666
$ws.Range("M115").Value = -479.6666
$ws.Range("H116").Value = 7378.8335
$ws.Range("J116").Value = 4499.5
$ws.Range("L116").Value = 4499.5
$ws.Range("N116").Value = -11383.5
$ws.Range("H135").Value = 2671.8333
$ws.Range("I135").Value = 2539.6
$ws.Range("J135").Value = 3333
$ws.Range("K135").Value = 22856.4
$ws.Range("L135").Value = 29997
$ws.Range("M135").Value = -20321.4
$ws.Range("N135").Value = -35067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6373.1035
$ws.Range("I32").Value = 4146.923
$ws.Range("J32").Value = 25666.666
$ws.Range("K32").Value = 4146.923
$ws.Range("L32").Value = 25666.666
$ws.Range("M32").Value = -3859.923
$ws.Range("N32").Value = -26240.666
$ws.Range("H61").Value = 2216.4443
$ws.Range("I61").Value = 2279.7144
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 2279.7144
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -2067.7144
$ws.Range("N61").Value = -2419
$ws.Range("H98").Value = 10000
$ws.Range("J98").Value = 10000
$ws.Range("L98").Value = 10000
$ws.Range("N98").Value = -15990
$ws.Range("H136").Value = 2216.4443
$ws.Range("I136").Value = 2279.7144
$ws.Range("J136").Value = 1995
$ws.Range("K136").Value = 6839.1432
$ws.Range("L136").Value = 5985
$ws.Range("M136").Value = -4289.1432
$ws.Range("N136").Value = -11085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 612.25
$ws.Range("I2").Value = 378.25
$ws.Range("J2").Value = 846.25
$ws.Range("K2").Value = 378.25
$ws.Range("L2").Value = 846.25
$ws.Range("M2").Value = -265.25
$ws.Range("N2").Value = -1072.25
$ws.Range("H31").Value = 6482.676
$ws.Range("I31").Value = 3825.7273
$ws.Range("J31").Value = 7606.769
$ws.Range("K31").Value = 3825.7273
$ws.Range("L31").Value = 7606.769
$ws.Range("M31").Value = -3530.7273
$ws.Range("N31").Value = -8196.769
$ws.Range("H34").Value = 6482.676
$ws.Range("I34").Value = 3825.7273
$ws.Range("J34").Value = 7606.769
$ws.Range("K34").Value = 3825.7273
$ws.Range("L34").Value = 7606.769
$ws.Range("M34").Value = -3623.7273
$ws.Range("N34").Value = -8010.769
$ws.Range("H39").Value = 1975
$ws.Range("I39").Value = 1975
$ws.Range("K39").Value = 1975
$ws.Range("M39").Value = -1584
$ws.Range("H49").Value = 1975
$ws.Range("I49").Value = 1975
$ws.Range("K49").Value = 1975
$ws.Range("M49").Value = -1793
$ws.Range("H58").Value = 4331
$ws.Range("I58").Value = 1999.1666
$ws.Range("K58").Value = 1999.1666
$ws.Range("M58").Value = -1796.1666
$ws.Range("H136").Value = 4331
$ws.Range("I136").Value = 1999.1666
$ws.Range("K136").Value = 5997.4998
$ws.Range("M136").Value = -3447.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64.59999999999999
$ws.Range("J2").Value = 35.833332
$ws.Range("L2").Value = 214.999992
$ws.Range("N2").Value = -440.999992
$ws.Range("H128").Value = 589994.2
$ws.Range("I128").Value = 589994.2
$ws.Range("K128").Value = 1769982.6
$ws.Range("M128").Value = -1765002.6
$ws.Range("H137").Value = 5058
$ws.Range("J137").Value = 5058
$ws.Range("L137").Value = 15174
$ws.Range("N137").Value = -25374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 178.5
$ws.Range("I16").Value = 178.5
$ws.Range("K16").Value = 178.5
$ws.Range("M16").Value = -8.5
$ws.Range("H55").Value = 1310.75
$ws.Range("I55").Value = 1818
$ws.Range("J55").Value = 465.33334
$ws.Range("K55").Value = 1818
$ws.Range("L55").Value = 465.33334
$ws.Range("M55").Value = -1645
$ws.Range("N55").Value = -811.33334
$ws.Range("H122").Value = 2983.7144
$ws.Range("I122").Value = 2925.5
$ws.Range("K122").Value = 8776.5
$ws.Range("M122").Value = -6326.5
$ws.Range("H132").Value = 4199
$ws.Range("I132").Value = 4998
$ws.Range("K132").Value = 14994
$ws.Range("M132").Value = -12464
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1950
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9312.5
$ws.Range("I62").Value = 4833.3335
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 4833.3335
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -4209.3335
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 9312.5
$ws.Range("I65").Value = 4833.3335
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 24166.6675
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -21046.6675
$ws.Range("N65").Value = -66240
$ws.Range("H96").Value = 1486.6666
$ws.Range("I96").Value = 1527.8572
$ws.Range("J96").Value = 1429
$ws.Range("K96").Value = 1527.8572
$ws.Range("L96").Value = 1429
$ws.Range("M96").Value = -154.8571999999999
$ws.Range("N96").Value = -4175
$ws.Range("H101").Value = 37800
$ws.Range("J101").Value = 37800
$ws.Range("L101").Value = 37800
$ws.Range("N101").Value = -44290
$ws.Range("H113").Value = 878.7143
$ws.Range("I113").Value = 1150.4
$ws.Range("K113").Value = 3451.2
$ws.Range("M113").Value = -1281.2
$ws.Range("H122").Value = 1499
$ws.Range("I122").Value = 1498.8334
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4496.5002
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2046.5002
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 3306.1
$ws.Range("I132").Value = 3117.889
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 9353.667000000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6823.667000000001
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2247.5715
$ws.Range("I136").Value = 2368.8
$ws.Range("J136").Value = 1944.5
$ws.Range("K136").Value = 7106.400000000001
$ws.Range("L136").Value = 5833.5
$ws.Range("M136").Value = -4556.400000000001
$ws.Range("N136").Value = -10933.5
$ws.Range("H139").Value = 59633.332
$ws.Range("J139").Value = 59633.332
$ws.Range("L139").Value = 59633.332
$ws.Range("N139").Value = -69913.33199999999
